# Weekly update: insert the latest week's "Poroto verde" price record at the
# top of the data block (row 105), pushing the existing rows 105-117 down to
# 106-118.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 105; existing rows shift down by one.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row with the latest week's data.
$ws.Range("A105").Value() = 8
$ws.Range("B105").Value() = "Terminal La Palmera de La Serena"
$ws.Range("C105").Value() = "Coquimbo"
$ws.Range("D105").Value() = 44449
$ws.Range("E105").Value() = 4
$ws.Range("F105").Value() = 100112031
$ws.Range("G105").Value() = "Poroto verde"
$ws.Range("H105").Value() = "Magnum"
$ws.Range("I105").Value() = "Primera"
$ws.Range("J105").Value() = 500
$ws.Range("K105").Value() = 34000
$ws.Range("L105").Value() = 35000
$ws.Range("M105").Value() = 34500
$ws.Range("N105").Value() = "$/malla 25 kilos"
$ws.Range("O105").Value() = "Perú"
$ws.Range("P105").Value() = 1380
$ws.Range("Q105").Value() = 25
$ws.Range("R105").Value() = "Hortaliza"
